$d = $word.ActiveDocument
$d.Content.Find.Execute(": Daty kampanii używające Gwiazdozbiór Bliźniąt 2022: 14-23 lutego, 14-24 marca", $true, $false, $false, $false, $false, $true, 1, $false, "2022: Daty kampanii używające Gwiazdozbiór Bliźniąt: 14-23 lutego, 14-24 marca", 2)
